$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Schneider: increase timeout to scrape result pages to 30s
# Update Next_update (column D) for Schneider-related rows (4,5,6)
# to match Last_update (column C), i.e. set to 44902 (2022-12-07)
$ws.Range("D4").Value = 44902
$ws.Range("D5").Value = 44902
$ws.Range("D6").Value = 44902

# Widen column C to fit content (closest achievable width to 25.1640625
# given this engine's pixel-quantized ColumnWidth storage)
$ws.Columns.Item(3).ColumnWidth = 24.33

# Move active selection to C14 as in the final file
$ws.Range("C14").Select()
